# Sample Project rule workbook: B11 ("R40") becomes the text "1".
# Because the target value "1" looks numeric, a plain .Value assignment
# would be auto-typed as a number by Excel and would also leave the
# original shared-string slot referencing "R40" pointed at a *number*
# cell. To faithfully reproduce a *text* cell (t="s") while leaving the
# cell's existing style (border/fill) completely untouched, stage the
# quoted text in an out-of-the-way helper cell, copy just its value
# (PasteSpecial -> values only) onto B11, then scrub the helper cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$helper = $ws.Cells.Item(500, 500)
$helper.Formula = "'1"          # leading apostrophe forces text storage
$helper.Copy()

$target = $ws.Cells.Item(11, 2) # B11
$target.PasteSpecial(-4163)     # xlPasteValues: value/type only, keep B11's own style

$excel.CutCopyMode = $false
$helper.Clear()
